# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F2").Value = 276
        $ws.Range("F3").Value = 170
        $ws.Range("F4").Value = 2030
        $ws.Range("F5").Value = 1636
        $ws.Range("F7").Value = 81
        $ws.Range("F8").Value = 645
        $ws.Range("F9").Value = 150
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F2").Value = 276
        $ws.Range("F3").Value = 170
        $ws.Range("F4").Value = 2030
        $ws.Range("F5").Value = 1636
        $ws.Range("F8").Value = 81
        $ws.Range("F9").Value = 645
        $ws.Range("F10").Value = 150
    }
}
